$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 27 de Marzo de 2020 a las 11:42'
$ws.Range("B4").Value = 85749
$ws.Range("C4").Value = 314
$ws.Range("E4").Value = 82577
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 1304
$ws.Range("B11").Value = 11951
$ws.Range("C11").Value = 140
$ws.Range("E11").Value = 11623
$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 197
$ws.Range("D15").Value = 858
$ws.Range("E15").Value = 6137
$ws.Range("F15").Value = 690
$ws.Range("B20").Value = 3423
$ws.Range("C20").Value = 51
$ws.Range("E20").Value = 3402
$ws.Range("B24").Value = 2858
$ws.Range("C24").Value = 18
$ws.Range("E24").Value = 2765
$ws.Range("A38").Value = 'Finlandia'
$ws.Range("B38").Value = 1035
$ws.Range("C38").Value = 77
$ws.Range("D38").Value = 10
$ws.Range("E38").Value = 1020
$ws.Range("F38").Value = 24
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 5
$ws.Range("A39").Value = 'Rumania'
$ws.Range("B39").Value = 1029
$ws.Range("D39").Value = 94
$ws.Range("E39").Value = 911
$ws.Range("F39").Value = 29
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 24
$ws.Range("A40").Value = 'Arabia Saudita'
$ws.Range("B40").Value = 1012
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 33
$ws.Range("E40").Value = 976
$ws.Range("F40").Value = 6
$ws.Range("H40").Value = 3
$ws.Range("B45").Value = 775
$ws.Range("C45").Value = 48
$ws.Range("E45").Value = 684
$ws.Range("A49").Value = 'Eslovenia'
$ws.Range("B49").Value = 632
$ws.Range("C49").Value = 70
$ws.Range("D49").Value = 10
$ws.Range("E49").Value = 616
$ws.Range("F49").Value = 14
$ws.Range("H49").Value = 6
$ws.Range("A50").Value = 'Argentina'
$ws.Range("B50").Value = 589
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 72
$ws.Range("E50").Value = 505
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 12
$ws.Range("A51").Value = 'Mexico'
$ws.Range("B51").Value = 585
$ws.Range("C51").Value = 110
$ws.Range("D51").Value = 4
$ws.Range("E51").Value = 573
$ws.Range("F51").Value = 1
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 8
$ws.Range("A52").Value = 'Peru'
$ws.Range("B52").Value = 580
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 14
$ws.Range("E52").Value = 557
$ws.Range("F52").Value = 14
$ws.Range("H52").Value = 9
$ws.Range("A53").Value = 'Estonia'
$ws.Range("B53").Value = 575
$ws.Range("C53").Value = 37
$ws.Range("D53").Value = 11
$ws.Range("E53").Value = 563
$ws.Range("F53").Value = 6
$ws.Range("H53").Value = 1
$ws.Range("A73").Value = 'Eslovaquia'
$ws.Range("B73").Value = 269
$ws.Range("C73").Value = 43
$ws.Range("D73").Value = 2
$ws.Range("E73").Value = 267
$ws.Range("F73").Value = 2
$ws.Range("H73").Value = 0
$ws.Range("A74").Value = 'Taiwan'
$ws.Range("B74").Value = 267
$ws.Range("C74").Value = 15
$ws.Range("D74").Value = 30
$ws.Range("E74").Value = 235
$ws.Range("F74").Value = 0
$ws.Range("H74").Value = 2
$ws.Range("A75").Value = 'Uruguay'
$ws.Range("B75").Value = 238
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 238
$ws.Range("F75").Value = 3
$ws.Range("H75").Value = 0
$ws.Range("A76").Value = 'Costa Rica'
$ws.Range("B76").Value = 231
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 2
$ws.Range("E76").Value = 227
$ws.Range("F76").Value = 5
$ws.Range("H76").Value = 2
$ws.Range("A77").Value = 'Bosnia y Herzegovina'
$ws.Range("B77").Value = 230
$ws.Range("C77").Value = 39
$ws.Range("D77").Value = 5
$ws.Range("E77").Value = 222
$ws.Range("F77").Value = 1
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 3
$ws.Range("A78").Value = 'Tunez'
$ws.Range("B78").Value = 227
$ws.Range("C78").Value = 30
$ws.Range("E78").Value = 219
$ws.Range("F78").Value = 10
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 6
$ws.Range("D81").Value = 5
$ws.Range("E81").Value = 208
$ws.Range("A132").Value = 'Isla de Man'
$ws.Range("B132").Value = 29
$ws.Range("C132").Value = 3
$ws.Range("D132").Value = 0
$ws.Range("E132").Value = 29
$ws.Range("A133").Value = 'Aruba'
$ws.Range("D133").Value = 1
$ws.Range("E133").Value = 27
$ws.Range("A134").Value = 'Guayana Francesa'
$ws.Range("B134").Value = 28
$ws.Range("D134").Value = 6
$ws.Range("E134").Value = 22
$ws.Range("A146").Value = 'El Salvador'
$ws.Range("A147").Value = 'Tanzania'
$ws.Range("A151").Value = 'Mongolia'
$ws.Range("A153").Value = 'Dominica'
$ws.Range("A158").Value = 'Guinea'
$ws.Range("C158").Value = 4
$ws.Range("A159").Value = 'Surinam'
$ws.Range("E159").Value = 8
$ws.Range("H159").Value = 0
$ws.Range("A160").Value = 'Islas Caimanes'
$ws.Range("D160").Value = 0
$ws.Range("E160").Value = 7
$ws.Range("H160").Value = 1
$ws.Range("A161").Value = 'Namibia'
$ws.Range("B161").Value = 8
$ws.Range("D161").Value = 2
$ws.Range("E161").Value = 6
$ws.Range("A162").Value = 'Mozambique'
$ws.Range("A163").Value = 'Antigua y Barbuda'
$ws.Range("A165").Value = 'Seychelles'
$ws.Range("E165").Value = 7
$ws.Range("H165").Value = 0
$ws.Range("A166").Value = 'Gabon'
$ws.Range("D166").Value = 0
$ws.Range("E166").Value = 6
$ws.Range("A167").Value = 'Curazao'
$ws.Range("B167").Value = 7
$ws.Range("D167").Value = 2
$ws.Range("E167").Value = 4
$ws.Range("H167").Value = 1
$ws.Range("A168").Value = 'Benin'
$ws.Range("A169").Value = 'Laos'
$ws.Range("A171").Value = 'Suazilandia'
$ws.Range("B171").Value = 6
$ws.Range("E171").Value = 6
$ws.Range("A172").Value = 'Fiyi'
$ws.Range("A174").Value = 'Birmania'
$ws.Range("A175").Value = 'Siria'
$ws.Range("E175").Value = 5
$ws.Range("H175").Value = 0
$ws.Range("A178").Value = 'Cabo Verde'
$ws.Range("B178").Value = 5
$ws.Range("H178").Value = 1
$ws.Range("A179").Value = 'Mali'
$ws.Range("A180").Value = 'Angola'
$ws.Range("A182").Value = 'Santa Sede'
$ws.Range("A184").Value = 'Republica del Chad'
$ws.Range("A185").Value = 'Somalia'
$ws.Range("C185").Value = 1
$ws.Range("A186").Value = 'Mauritania'
$ws.Range("A188").Value = 'Butan'
$ws.Range("C188").Value = 1
$ws.Range("A189").Value = 'San Martin (Parte Holandesa)'
$ws.Range("C189").Value = 0
$ws.Range("A190").Value = 'Republica de Africa Central'
$ws.Range("C190").Value = 0
$ws.Range("A191").Value = 'Gambia'
$ws.Range("D191").Value = 0
$ws.Range("H191").Value = 1
$ws.Range("A192").Value = 'Sudan'
$ws.Range("A193").Value = 'Nepal'
$ws.Range("A194").Value = 'Santa Lucia'
$ws.Range("D194").Value = 1
$ws.Range("H194").Value = 0
$ws.Range("A195").Value = 'Belice'
$ws.Range("A196").Value = 'Islas Virgenes Britanicas'
$ws.Range("A197").Value = 'Guinea-Bisau'
$ws.Range("A198").Value = 'San Cristobal y Nieves'
$ws.Range("A199").Value = 'Islas Turcas y Caicos'
$ws.Range("A200").Value = 'Anguila'
$ws.Range("A202").Value = 'Papua Nueva Guinea'
$ws.Range("A203").Value = 'Libia'
$ws.Range("A204").Value = 'Timor Oriental'
$ws.Range("A205").Value = 'San Vicente y las Granadinas'
